$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: some updated cell values look like plain numbers (e.g. "593.80",
# "0.0000240") even though the column stores text. Excel's COM layer will
# silently convert such strings to numeric values on assignment, which would
# lose the original text formatting (trailing zeros, exact digit count, etc).
# To avoid that we temporarily force the cell to Text format, assign the
# string, then restore the default "Normal" style so the cell ends up with
# no explicit style (matching the source workbook) but keeps the literal text.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '67.145.39'
$ws.Range("E2").Value = '  -3.58%  '

$ws.Range("D3").Value = '3.676.80'
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  -0.02%  '

Set-TextValue $ws.Range("D5") '593.80'
$ws.Range("E5").Value = '  -3.27%  '

Set-TextValue $ws.Range("D6") '165.50'
$ws.Range("E6").Value = '  -6.65%  '

$ws.Range("D7").Value = '3.677.82'
$ws.Range("E7").Value = '  -3.36%  '

$ws.Range("E8").Value = '  -0.09%  '

Set-TextValue $ws.Range("D9") '0.524'
$ws.Range("E9").Value = '  -1.02%  '

Set-TextValue $ws.Range("D10") '0.159'
$ws.Range("E10").Value = '  -4.96%  '

Set-TextValue $ws.Range("D11") '6.14'
$ws.Range("E11").Value = '  -5.35%  '

$ws.Range("E12").Value = '  -4.85%  '

Set-TextValue $ws.Range("D13") '37.44'
$ws.Range("E13").Value = '  -5.78%  '

Set-TextValue $ws.Range("D14") '0.0000240'
$ws.Range("E14").Value = '  -5.92%  '

$ws.Range("D15").Value = '4.284.78'
$ws.Range("E15").Value = '  -3.39%  '

$ws.Range("D16").Value = '3.669.45'
$ws.Range("E16").Value = '  -3.48%  '

$ws.Range("D17").Value = '67.161.25'
$ws.Range("E17").Value = '  -3.67%  '

Set-TextValue $ws.Range("D18") '7.14'
$ws.Range("E18").Value = '  -5.49%  '

$ws.Range("E19").Value = '  -4.11%  '

Set-TextValue $ws.Range("D20") '17.12'
$ws.Range("E20").Value = '  +3.03%  '

Set-TextValue $ws.Range("D21") '489.38'
$ws.Range("E21").Value = '  -3.42%  '

Set-TextValue $ws.Range("D22") '9.10'
$ws.Range("E22").Value = '  -5.45%  '

Set-TextValue $ws.Range("D23") '0.715'
$ws.Range("E23").Value = '  -2.79%  '

Set-TextValue $ws.Range("D24") '85.40'
$ws.Range("E24").Value = '  -1.04%  '

Set-TextValue $ws.Range("D25") '2.29'
$ws.Range("E25").Value = '  -7.15%  '

Set-TextValue $ws.Range("D26") '0.0000138'
$ws.Range("E26").Value = '  -4.57%  '

Set-TextValue $ws.Range("D27") '12.09'
$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("E28").Value = '  -0.43%  '

Set-TextValue $ws.Range("D29") '9.92'
$ws.Range("E29").Value = '  -6.03%  '

Set-TextValue $ws.Range("D30") '2.91'
$ws.Range("E30").Value = '  -2.30%  '

Set-TextValue $ws.Range("D31") '2.36'
$ws.Range("E31").Value = '  -6.45%  '

$ws.Range("E32").Value = '  -4.24%  '

Set-TextValue $ws.Range("D33") '31.61'
$ws.Range("E33").Value = '  +0.39%  '

$ws.Range("D34").Value = '3.807.03'
$ws.Range("E34").Value = '  -3.54%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D35") '0.106'
$ws.Range("E35").Value = '  -6.58%  '

$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.607.35'
$ws.Range("E36").Value = '  -3.52%  '

Set-TextValue $ws.Range("D37") '0.999'
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("E38").Value = '  -5.60%  '

$ws.Range("E39").Value = '  -6.16%  '

$ws.Range("E40").Value = '  -6.99%  '

Set-TextValue $ws.Range("D41") '0.322'
$ws.Range("E41").Value = '  -4.66%  '

Set-TextValue $ws.Range("D42") '437.32'
$ws.Range("E42").Value = '  -9.16%  '

Set-TextValue $ws.Range("D43") '48.59'
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("E44").Value = '  -6.82%  '

Set-TextValue $ws.Range("D45") '2.77'
$ws.Range("E45").Value = '  -8.35%  '

Set-TextValue $ws.Range("D46") '8.32'
$ws.Range("E46").Value = '  -2.84%  '

$ws.Range("E47").Value = '  +0.00%  '

Set-TextValue $ws.Range("D48") '142.39'
$ws.Range("E48").Value = '  +1.98%  '

Set-TextValue $ws.Range("D49") '39.72'
$ws.Range("E49").Value = '  -9.85%  '

$ws.Range("D50").Value = '2.748.96'
$ws.Range("E50").Value = '  -6.15%  '

$ws.Range("E51").Value = '  -4.64%  '
